# Updated cryptos list values (Price / Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking strings ("1.011", "0.07230", "8.250", ...).
# Setting NumberFormat to Text ("@") before writing keeps them as text (as in the
# original inline-string cells) instead of Excel auto-converting to a Number and
# silently dropping meaningful trailing zeros / thousands-style dots.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.991.53'
$ws.Range("E2").Value = '  +1.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.06'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.011'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '309.27'
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4777'
$ws.Range("E7").Value = '  +2.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3675'
$ws.Range("E8").Value = '  +2.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07230'
$ws.Range("E9").Value = '  +1.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9319'
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.76'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07735'
$ws.Range("E12").Value = '  +1.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.831.49'
$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.347'
$ws.Range("E14").Value = '  +1.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.437'
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.01'
$ws.Range("E16").Value = '  +1.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.013'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008645'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("E19").Value = '  +0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.017.07'
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("E21").Value = '  +1.79%  '

$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.64'
$ws.Range("E23").Value = '  +0.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.924'
$ws.Range("E24").Value = '  +0.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.03'
$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.22'
$ws.Range("E26").Value = '  +1.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.019'
$ws.Range("E27").Value = '  +1.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.38'
$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08859'
$ws.Range("E30").Value = '  +0.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.314'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.184'
$ws.Range("E32").Value = '  +1.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7418'
$ws.Range("E33").Value = '  +0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.502'
$ws.Range("E34").Value = '  +1.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.736'
$ws.Range("E35").Value = '  -4.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.111'
$ws.Range("E36").Value = '  +3.18%  '

$ws.Range("E37").Value = '  +1.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05262'
$ws.Range("E38").Value = '  +2.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.975'
$ws.Range("E39").Value = '  +0.86%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5261'
$ws.Range("E40").Value = '  +3.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.015'
$ws.Range("E41").Value = '  +1.58%  '

$ws.Range("E42").Value = '  +1.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.250'
$ws.Range("E43").Value = '  +1.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.62'
$ws.Range("E44").Value = '  +5.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4753'
$ws.Range("E45").Value = '  +2.15%  '

$ws.Range("E46").Value = '  +0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.71'
$ws.Range("E47").Value = '  +3.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.609'
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.55'
$ws.Range("E49").Value = '  +2.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06071'
$ws.Range("E50").Value = '  +0.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8902'
$ws.Range("E51").Value = '  +4.01%  '
